$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add new header cells (copy header format from an existing header cell first) ---
$ws.Range("G1").Copy()
$ws.Range("H1:K1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- 2. Update header text ---
$ws.Range("F1").Value = "Link Rank"
$ws.Range("G1").Value = "Duration Rank"
$ws.Range("H1").Value = "Views Rank"
$ws.Range("I1").Value = "Link First"
$ws.Range("J1").Value = "Duration First"
$ws.Range("K1").Value = "Views First"

# --- 3. Shrink font across the table (header + data) to size 10 ---
$ws.Range("A1:K1").Font.Size = 10
$ws.Range("A2:E4").Font.Size = 10

# --- 4. Column widths for the new/changed columns ---
# (target stored widths: 41.5703125 / 16.28515625 / 14.140625 -- the runtime
# quantizes ColumnWidth to 1/6-character steps, so these inputs land on the
# closest achievable stored width.)
$ws.Columns.Item(6).ColumnWidth = 40.666666666666664
$ws.Columns.Item(7).ColumnWidth = 15.5
$ws.Columns.Item(8).ColumnWidth = 13.333333333333334
$ws.Columns.Item(9).ColumnWidth = 40.666666666666664
$ws.Columns.Item(10).ColumnWidth = 15.5
$ws.Columns.Item(11).ColumnWidth = 13.333333333333334

# --- 5. Page setup orientation ---
$ws.PageSetup.Orientation = 1

# --- 6. Selection cosmetic change ---
$ws.Range("K2").Select()
